# Mise à jour de l'application
# Adds a new attendance-day column (BQ) for 2025-10-24 (serial 45954),
# one column to the right of the existing last day column (BP, 2025-10-22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Attendance mark for each player row in the new BQ column.
# Row 12 (Yanis Berrached) is intentionally omitted: that player's data
# already stopped before column BP, so no entry is added for them either.
$attendance = [ordered]@{
    2  = "P"   # Alban Rambaud
    3  = "P"   # Jassim Assoul
    4  = "P"   # Enzo Vita
    5  = "P"   # Romain Thunet
    6  = "B"   # Amine Taiar
    7  = "P"   # Naim Ighbane
    8  = "B"   # Hedi Nasri
    9  = "P"   # Mattheo Haon
    10 = "P"   # Maé Clavel
    11 = "P"   # Levy Ndoutoume
    13 = "B"   # Rayane Chayebi
    14 = "P"   # Ilan Ihaddadene
    15 = "P"   # Karahali Souaré
    16 = "P"   # Amir Etien
    17 = "RH"  # Karim Belmahi
    18 = "P"   # Emmanuel Valey
    19 = "P"   # Jeremie Laurent
    20 = "P"   # Sofiane Belle
    21 = "B"   # Amir Kherrab
    22 = "P"   # Naim Dhib
    23 = "B"   # Wael Fareh
    24 = "B"   # Yoan Zouma
    25 = "P"   # Ilyes Boughanmi
    26 = "P"   # Omar Benyounes
    27 = "P"   # Yoann Martelat
    28 = "P"   # Malik Boussaid
    29 = "P"   # Kamal Bafounta
}

# 1) Write every value first (direct .Value assignment keeps the wide
#    COUNTA/COUNTIF formulas in columns B:J dirty so they recalc correctly).
$ws.Range("BQ1").Value = 45954
foreach ($row in $attendance.Keys) {
    $ws.Range("BQ$row").Value = $attendance[$row]
}

# 2) Now copy number-format/alignment from the previous day's column (BP)
#    onto the new column (BQ), cell by cell, so each new cell reuses the
#    same style as its BP counterpart instead of creating new styles.
$ws.Range("BP1").Copy() | Out-Null
$ws.Range("BQ1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

foreach ($row in $attendance.Keys) {
    $ws.Range("BP$row").Copy() | Out-Null
    $ws.Range("BQ$row").PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

# 3) Match the author's final selection/cursor position.
$ws.Range("BT25").Select() | Out-Null
